$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting existing rows 155.. down by one.
$ws.Rows.Item(155).Insert()

# Populate the new row 155 with the new record's data.
$ws.Cells.Item(155, 1).Value = 3
$ws.Cells.Item(155, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44529
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 100112043
$ws.Cells.Item(155, 7).Value = "Pepino ensalada"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 125
$ws.Cells.Item(155, 11).Value = 6000
$ws.Cells.Item(155, 12).Value = 6500
$ws.Cells.Item(155, 13).Value = 6260
$ws.Cells.Item(155, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 89
$ws.Cells.Item(155, 17).Value = 70
$ws.Cells.Item(155, 18).Value = "Hortaliza"
